$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink on C2 (and its relationship) before rewriting
# the sheet contents.
$ws.Hyperlinks.Delete()

# Clear the columns (F:I) that are being dropped from the sheet so the used
# range / dimension shrinks back down to A1:E2.
$ws.Range("F1:I2").ClearContents()

# --- Row 1 (headers) ---
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "locacalizacion"
$ws.Range("C1").Value = "Correo electrónico"
$ws.Range("D1").Value = "id"
$ws.Range("E1").Value = "kind"

# --- Row 2 (data) ---
$ws.Range("A2").ClearContents()
$ws.Range("B2").Value = "18:13:14:12S"
$ws.Range("C2").Value = "jorge@email.es"
$ws.Range("D2").Value = "ID4"
$ws.Range("E2").Value = 1
